$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Tested?" column (B) for the gripper / elevator related rows
$ws.Range("B14").Value = "yes"
$ws.Range("B15").Value = "yes"
$ws.Range("B16").Value = "yes"
$ws.Range("B17").Value = "no"
$ws.Range("B19").Value = "yes"
$ws.Range("B20").Value = "no"
$ws.Range("B21").Value = "yes"
$ws.Range("B22").Value = "no"
$ws.Range("B23").Value = "no"
$ws.Range("B24").Value = "no"

# Apply the built-in "Comma" style to B14 (matches the Comma cell style added to styles.xml)
$ws.Range("B14").Style = "Comma"

# Update the view: move/scroll the selection to B26
$ws.Activate()
$ws.Range("B26").Select()
